$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Toggle Grupo_Experimental (column B) values for several rows
$ws.Range("B2").Value = "Con SmartScore"
$ws.Range("B3").Value = "Sin SmartScore"
$ws.Range("B4").Value = "Sin SmartScore"
$ws.Range("B7").Value = "Sin SmartScore"
$ws.Range("B9").Value = "Con SmartScore"
$ws.Range("B10").Value = "Con SmartScore"
$ws.Range("B11").Value = "Sin SmartScore"

# Convert the SmartScore text values in row 11 to true numbers
$ws.Range("I11").Value = 0.579
$ws.Range("L11").Value = 0.479
$ws.Range("O11").Value = 0.469
$ws.Range("R11").Value = 0.601
$ws.Range("U11").Value = 0.5590000000000001
$ws.Range("X11").Value = 0.547
$ws.Range("AD11").Value = 0.59
$ws.Range("AG11").Value = 0.5669999999999999
